# "Lock in current version" of the 18_cues cue-sequence sheet.
#
# The experiment-generation step re-drew a fresh, shuffled sequence of
# word / image / category cues. This script writes that finalized sequence
# into the 49-row table (row 1 is the "word" | "image" | "category" header
# and is left untouched; data lives in rows 2-49).
#
# The new values are written one whole column at a time (all of column A,
# then all of column B, then all of column C) so that any brand-new text
# introduced by this edit is appended to the workbook's shared-string table
# in the same left-to-right, top-to-bottom order the source data was drawn in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: word
$words = @(
    "wenden", "quälen", "seufzen", "rufen", "parken", "beißen",
    "enden", "schneiden", "stillen", "opfern", "schreiben", "triefen",
    "schalten", "spenden", "fühlen", "klagen", "decken", "albern",
    "sparen", "packen", "reizen", "hören", "regnen", "kümmern",
    "drohen", "tollen", "deuten", "orten", "rasen", "achten",
    "weigern", "zählen", "bellen", "ächzen", "herrschen", "tragen",
    "kosten", "gleichen", "knien", "dauern", "holen", "münzen",
    "stören", "erben", "wehtun", "bremsen", "schleppen", "mühen"
)

# Column B: image
$images = @(
    "none", "face/face020.jpg", "house/house003.jpg", "none", "face/face003.jpg", "face/face026.jpg",
    "none", "face/face011.jpg", "face/face015.jpg", "none", "face/face014.jpg", "face/face024.jpg",
    "none", "house/house008.jpg", "face/face007.jpg", "none", "house/house017.jpg", "face/face021.jpg",
    "none", "face/face028.jpg", "house/house031.jpg", "none", "house/house026.jpg", "house/house020.jpg",
    "none", "house/house021.jpg", "face/face017.jpg", "none", "face/face005.jpg", "house/house012.jpg",
    "none", "face/face009.jpg", "house/house010.jpg", "none", "house/house018.jpg", "house/house023.jpg",
    "none", "house/house014.jpg", "face/face029.jpg", "none", "face/face000.jpg", "face/face002.jpg",
    "none", "house/house024.jpg", "house/house007.jpg", "none", "house/house028.jpg", "house/house022.jpg"
)

# Column C: category
$categories = @(
    "none", "face", "house", "none", "face", "face",
    "none", "face", "face", "none", "face", "face",
    "none", "house", "face", "none", "house", "face",
    "none", "face", "house", "none", "house", "house",
    "none", "house", "face", "none", "face", "house",
    "none", "face", "house", "none", "house", "house",
    "none", "house", "face", "none", "face", "face",
    "none", "house", "house", "none", "house", "house"
)

$firstDataRow = 2

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + $firstDataRow, 1).Value = $words[$i]
}
for ($i = 0; $i -lt $images.Length; $i++) {
    $ws.Cells.Item($i + $firstDataRow, 2).Value = $images[$i]
}
for ($i = 0; $i -lt $categories.Length; $i++) {
    $ws.Cells.Item($i + $firstDataRow, 3).Value = $categories[$i]
}
